$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the whole data block (A1:E4) down by one row and right by one column
# (i.e. it now occupies B2:F5), leaving row 1 / column A empty. Inserting a
# whole column/row shifts the existing cells (values, types and styles)
# along with it, rather than re-creating the content from scratch.
$ws.Columns("A").Insert() | Out-Null
$ws.Rows("1").Insert() | Out-Null

$ws.Range("B2").Select() | Out-Null
